$wb = $excel.ActiveWorkbook

$wsWater = $wb.Worksheets.Item("watercolours")
$wsPaints = $wb.Worksheets.Item("paints")

# Rename "Violet Purple" -> "Eggplant Purple" and add new "Egg Blue" demo variable
$wsWater.Range("C4").Value = "Egg Blue"
$wsWater.Range("B3").Value = "Eggplant Purple"

# Update selection on paints sheet (no longer the active tab)
$wsPaints.Range("E25").Select()

# Make watercolours the active sheet/tab, select B3
$wsWater.Activate()
$wsWater.Range("B3").Select()
